$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211, shifting existing rows 211:325 down to 212:326
$ws.Rows(211).Insert()

# Populate the newly inserted row 211 with the new data record
$ws.Cells.Item(211, 1).Value = 7
$ws.Cells.Item(211, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(211, 3).Value = "Ñuble"
$ws.Cells.Item(211, 4).Value = 44813
$ws.Cells.Item(211, 5).Value = 16
$ws.Cells.Item(211, 6).Value = 100114013
$ws.Cells.Item(211, 7).Value = "Zanahoria"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 120
$ws.Cells.Item(211, 11).Value = 9000
$ws.Cells.Item(211, 12).Value = 10000
$ws.Cells.Item(211, 13).Value = 9500
$ws.Cells.Item(211, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(211, 15).Value = "Región de Ñuble"
$ws.Cells.Item(211, 16).Value = 475
$ws.Cells.Item(211, 17).Value = 20
$ws.Cells.Item(211, 18).Value = "Hortaliza"
